$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add headers for new columns I (I0) and J (IF), matching style of existing headers (bold, bordered, centered)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Values for column I (rows 2-32)
$iValues = @(6,9,7,6,8,8,9,2,8,1,9,5,8,7,5,1,1,7,7,9,7,1,1,1,1,4,9,9,7,7,1)
# Values for column J (rows 2-32)
$jValues = @(7,9,7,6,8,8,9,3,8,3,9,7,9,8,6,3,3,9,8,9,7,3,1,1,2,5,9,9,7,8,2)

for ($idx = 0; $idx -lt $iValues.Length; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$idx]
    $ws.Cells.Item($row, 10).Value = $jValues[$idx]
}
